$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the train schedule row 4 (y_corrSteps, y_nrSteps, alienID) to the
# corrected values from the task edit.
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 13

# The author's last click/selection ended up on E4 (was I7) before saving.
$ws.Range("E4").Select()
